# Dashboard redesign for the "bot brainbow" sheet:
#  - the first two data rows no longer carry an Id value in column A,
#    so clear out A2:A3 (they have no explicit style, so clearing removes
#    the cell entries entirely, matching the source edit).
#  - the sheet's remembered selection moves from A7 (off the used range)
#    to A3 (the last populated row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A3").ClearContents() | Out-Null

$ws.Range("A3").Select() | Out-Null
